$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Swap the contents of rows 195 and 196 (everything except the
#    Indice/pais/torneio/temporada/data_partida columns A:E, which stay
#    put). Row 195 used to be "Envigado vs Dep. Pasto", row 196 used to
#    be "Chico vs Dep. Cali" — the edit swaps their order in the sheet.
# ------------------------------------------------------------------

$envigadoMatch = @{
    F = "Envigado"
    H = "Dep. Pasto"
    J = 2.7
    K = "01/11/2023 10:42"
    L = 2.83
    M = "08/11/2023 01:29"
    N = 2.92
    O = "01/11/2023 10:42"
    P = 3.26
    Q = "08/11/2023 01:29"
    R = 3.03
    S = "01/11/2023 10:42"
    T = 2.69
    U = "08/11/2023 01:29"
    V = "https://www.betexplorer.com/football/colombia/primera-a/envigado-dep-pasto/ptUg3o6n/"
}

$chicoMatch = @{
    F = "Chico"
    H = "Dep. Cali"
    J = 2.65
    K = "01/11/2023 10:43"
    L = 3.75
    M = "08/11/2023 01:23"
    N = 2.97
    O = "01/11/2023 10:43"
    P = 3.1
    Q = "08/11/2023 01:23"
    R = 2.9
    S = "01/11/2023 10:43"
    T = 2.25
    U = "08/11/2023 01:23"
    V = "https://www.betexplorer.com/football/colombia/primera-a/chico-dep-cali/pGp3UN4o/"
}

function Set-MatchRow($ws, $rowNum, $match) {
    $ws.Range("F$rowNum").Value = $match.F
    $ws.Range("H$rowNum").Value = $match.H
    $ws.Range("J$rowNum").Value = $match.J
    $ws.Range("K$rowNum").Value = $match.K
    $ws.Range("L$rowNum").Value = $match.L
    $ws.Range("M$rowNum").Value = $match.M
    $ws.Range("N$rowNum").Value = $match.N
    $ws.Range("O$rowNum").Value = $match.O
    $ws.Range("P$rowNum").Value = $match.P
    $ws.Range("Q$rowNum").Value = $match.Q
    $ws.Range("R$rowNum").Value = $match.R
    $ws.Range("S$rowNum").Value = $match.S
    $ws.Range("T$rowNum").Value = $match.T
    $ws.Range("U$rowNum").Value = $match.U
    $ws.Range("V$rowNum").Value = $match.V
}

Set-MatchRow $ws 195 $chicoMatch
Set-MatchRow $ws 196 $envigadoMatch

# ------------------------------------------------------------------
# 2. Append two new match rows (210 and 211) after the current last
#    row (209), copying the row-209 cell formatting first so the new
#    rows pick up the same styles (Indice column style, date column
#    style) as the rest of the sheet.
# ------------------------------------------------------------------

$ws.Range("A209:V209").Copy() | Out-Null
$ws.Range("A210:V210").PasteSpecial(-4122) | Out-Null
$ws.Range("A211:V211").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

function Set-NewRow($ws, $rowNum, $data) {
    $ws.Range("A$rowNum").Value = $data.Indice
    $ws.Range("B$rowNum").Value = $data.Pais
    $ws.Range("C$rowNum").Value = $data.Torneio
    # "temporada" is stored as text ("2023"), matching every other row in
    # the sheet, even though it looks numeric — apply a text number format
    # first so Excel doesn't coerce the value to a number.
    $ws.Range("D$rowNum").NumberFormat = "@"
    $ws.Range("D$rowNum").Value = $data.Temporada
    $ws.Range("E$rowNum").Value = $data.DataPartida
    $ws.Range("F$rowNum").Value = $data.Home
    $ws.Range("G$rowNum").Value = $data.HomeGols
    $ws.Range("H$rowNum").Value = $data.Away
    $ws.Range("I$rowNum").Value = $data.AwayGols
    $ws.Range("J$rowNum").Value = $data.HomeOpeningOdds
    $ws.Range("K$rowNum").Value = $data.HomeOpeningDataHora
    $ws.Range("L$rowNum").Value = $data.HomeClosingOdds
    $ws.Range("M$rowNum").Value = $data.HomeClosingDataHora
    $ws.Range("N$rowNum").Value = $data.DrawOpeningOdds
    $ws.Range("O$rowNum").Value = $data.DrawOpeningDataHora
    $ws.Range("P$rowNum").Value = $data.DrawClosingOdds
    $ws.Range("Q$rowNum").Value = $data.DrawClosingDataHora
    $ws.Range("R$rowNum").Value = $data.AwayOpeningOdds
    $ws.Range("S$rowNum").Value = $data.AwayOpeningDataHora
    $ws.Range("T$rowNum").Value = $data.AwayClosingOdds
    $ws.Range("U$rowNum").Value = $data.AwayClosingDataHora
    $ws.Range("V$rowNum").Value = $data.UrlPartida
}

$row210 = @{
    Indice = 209
    Pais = "colombia"
    Torneio = "primera-a"
    Temporada = "2023"
    DataPartida = 45255.95833333334
    Home = "Deportes Tolima"
    HomeGols = 4
    Away = "Dep. Cali"
    AwayGols = 2
    HomeOpeningOdds = 1.82
    HomeOpeningDataHora = "19/11/2023 01:42"
    HomeClosingOdds = 1.65
    HomeClosingDataHora = "25/11/2023 22:55"
    DrawOpeningOdds = 3.47
    DrawOpeningDataHora = "19/11/2023 01:42"
    DrawClosingOdds = 3.94
    DrawClosingDataHora = "25/11/2023 22:55"
    AwayOpeningOdds = 4.8
    AwayOpeningDataHora = "19/11/2023 01:42"
    AwayClosingOdds = 5.61
    AwayClosingDataHora = "25/11/2023 22:55"
    UrlPartida = "https://www.betexplorer.com/football/colombia/primera-a/deportes-tolima-dep-cali/vyoEzd9h/"
}

$row211 = @{
    Indice = 210
    Pais = "colombia"
    Torneio = "primera-a"
    Temporada = "2023"
    DataPartida = 45256.0625
    Home = "Junior"
    HomeGols = 3
    Away = "Aguilas"
    AwayGols = 1
    HomeOpeningOdds = 1.95
    HomeOpeningDataHora = "19/11/2023 05:42"
    HomeClosingOdds = 1.89
    HomeClosingDataHora = "26/11/2023 01:20"
    DrawOpeningOdds = 3.27
    DrawOpeningDataHora = "19/11/2023 05:42"
    DrawClosingOdds = 3.34
    DrawClosingDataHora = "26/11/2023 01:19"
    AwayOpeningOdds = 4.42
    AwayOpeningDataHora = "19/11/2023 05:42"
    AwayClosingOdds = 4.82
    AwayClosingDataHora = "26/11/2023 01:20"
    UrlPartida = "https://www.betexplorer.com/football/colombia/primera-a/junior-aguilas-doradas/IRqIZFOb/"
}

Set-NewRow $ws 210 $row210
Set-NewRow $ws 211 $row211
